# Swap the order of names in the "Recorded By" column (G) so that
# "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com"
# across the whole used range of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Columns.Item(7)  # Column G - "Recorded By"
$range.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
